# Congreso_Resultados.xlsx update
# - lowercases header labels, adds trailing spaces to several narrative labels
# - replaces the evaluation metrics with a new set of results (accuracy run)
# - applies a 4-decimal number format to the metric/time columns
# - reworks the formatting/merges of the "donde:" legend block
# - updates the current selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text updates (shared strings)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "no."
$ws.Range("B1").Value = "algoritmo"

$ws.Range("B2").Value  = "Ridge regression and classification "
$ws.Range("B3").Value  = "Classification, Support Vector Machines "
$ws.Range("B4").Value  = "Classification, Stochastic Gradient Descent "
$ws.Range("B5").Value  = "Nearest Neighbours Classification "
$ws.Range("B6").Value  = "Gaussian Process Classification "
$ws.Range("B7").Value  = "Classification, Decision Trees "
$ws.Range("B8").Value  = "Voting Classifier (Ensemble) "
$ws.Range("B9").Value  = "Classification, Neural network models "

$ws.Range("B11").Value = "valores máximos "
$ws.Range("B12").Value = "donde: "
$ws.Range("B13").Value = "precision: precisión de las predicciones positivas "
$ws.Range("B14").Value = "recall: fracción de predicciones positivas correctamente identificadas "
$ws.Range("B15").Value = "f1-score: media armónica de precision y recall (F measure) "
$ws.Range("B16").Value = "tiempo: tiempo de proceso "

# ---------------------------------------------------------------------------
# 2. New metric values for each algorithm (rows 2-9) and the "valores
#    maximos" summary row (row 11)
# ---------------------------------------------------------------------------
$ws.Range("C2").Value2 = 0.8333 ; $ws.Range("D2").Value2 = 0.3333 ; $ws.Range("E2").Value2 = 0.4762 ; $ws.Range("F2").Value2 = 7.9737
$ws.Range("C3").Value2 = 1      ; $ws.Range("D3").Value2 = 0.08   ; $ws.Range("E3").Value2 = 0.1481 ; $ws.Range("F3").Value2 = 66.8335
$ws.Range("C4").Value2 = 1      ; $ws.Range("D4").Value2 = 0.5    ; $ws.Range("E4").Value2 = 0.6667 ; $ws.Range("F4").Value2 = 4.3036
$ws.Range("C5").Value2 = 0.2727 ; $ws.Range("D5").Value2 = 0.1429 ; $ws.Range("E5").Value2 = 0.1875 ; $ws.Range("F5").Value2 = 3.5404
$ws.Range("C6").Value2 = 0.8889 ; $ws.Range("D6").Value2 = 0.8889 ; $ws.Range("E6").Value2 = 0.8889 ; $ws.Range("F6").Value2 = 4.2267
$ws.Range("C7").Value2 = 0.4545 ; $ws.Range("D7").Value2 = 0.4545 ; $ws.Range("E7").Value2 = 0.4545 ; $ws.Range("F7").Value2 = 24.455
$ws.Range("C8").Value2 = 0.9    ; $ws.Range("D8").Value2 = 0.9534 ; $ws.Range("E8").Value2 = 0.9259 ; $ws.Range("F8").Value2 = 11.8827
$ws.Range("C9").Value2 = 0.9398 ; $ws.Range("D9").Value2 = 0.8797 ; $ws.Range("E9").Value2 = 0.9087 ; $ws.Range("F9").Value2 = 1.3256

$ws.Range("C11").Value2 = 1      ; $ws.Range("D11").Value2 = 0.9534 ; $ws.Range("E11").Value2 = 0.9259 ; $ws.Range("F11").Value2 = 66.8335

# apply the 4-decimal number format to every metric/time cell that now holds
# a value
$ws.Range("C2:F9").NumberFormat = "0.0000"
$ws.Range("C11:F11").NumberFormat = "0.0000"
$ws.Range("C10:F10").NumberFormat = "0.0000"

# ---------------------------------------------------------------------------
# 3. Blank separator row (row 10): drop the stray alignment flag on A10:B10
#    and give C10:F10 the same 4-decimal format as the data rows (done above)
# ---------------------------------------------------------------------------
$ws.Range("A10:B10").Font.Name = "Arial"

# ---------------------------------------------------------------------------
# 4. "donde:" legend block (rows 13-15): switch the affected cells to Times
#    New Roman and re-merge the label ranges
# ---------------------------------------------------------------------------
$ws.Range("A10:B10").Select() | Out-Null
$ws.Range("A10:B10").UnMerge()

$legend1 = $ws.Range("B13:C13")
$legend1.Font.Name = "Times New Roman"
$legend1.WrapText = $false
$legend1.Merge()

$legend2 = $ws.Range("B14:E14")
$legend2.Font.Name = "Times New Roman"
$legend2.WrapText = $false
$legend2.Merge()

$legend3 = $ws.Range("B15:D15")
$legend3.Font.Name = "Times New Roman"
$legend3.WrapText = $false
$legend3.Merge()

$ws.Range("B16").Font.Name = "Times New Roman"

# ---------------------------------------------------------------------------
# 5. Selection
# ---------------------------------------------------------------------------
$ws.Range("C2:F11").Select()
